$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Finally conclusion is discussed in Section 7." ->
#    "Finally conclusion is given in Section 7."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("lusion is discussed in Section 7", $true, $false, $false, $false, $false, `
    $true, 1, $false, "lusion is given in Section 7", 1) | Out-Null

# ---------------------------------------------------------------------
# 2) Anonymity bullet: "...identities to connection card issuer (TTP)..."
#    -> "...identities to Trusted Third Party (TTP)..." and the "TTP"
#    run right after the new "(" loses its italic formatting.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("connection card issuer", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Trusted Third Party", 1) | Out-Null

$scanRange = $d.Content
$scanRange.Find.Execute("identities to Trusted Third Party (", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
if ($scanRange.Find.Found) {
    $scanRange.Collapse(0)
    $ttpRange = $d.Range($scanRange.Start, $scanRange.Start + 3)
    if ($ttpRange.Text -eq "TTP") {
        $ttpRange.Italic = 0
    }
}

# ---------------------------------------------------------------------
# 3) Remove the empty "IEEEParagraph"-styled paragraph that sits right
#    before the "Conclusion" heading.
# ---------------------------------------------------------------------
$headingRange = $d.Content
$headingRange.Find.Execute("Conclusion", $true, $true, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
if ($headingRange.Find.Found) {
    $headingPara = $headingRange.Paragraphs.First
    $prevPara = $headingPara.Previous()
    if ($prevPara -ne $null -and $prevPara.Range.Text.Trim().Length -eq 0) {
        $prevPara.Range.Delete() | Out-Null
    }
}

# ---------------------------------------------------------------------
# 4) Merge the two runs forming the Wi-Fi Reports/Papers hyperlink text
#    into a single run (no visible text change, just a no-op edit that
#    coalesces the adjacent identically-formatted runs).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Wi-Fi_Reports_and_Papers", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Wi-Fi_Reports_and_Papers", 1) | Out-Null
